# Fuel mixing assumptions: add international bunkers (ship/air) rows for
# the 2040, 2050 and 2070 milestone years to the international_supply_side
# sheet, then refresh the sheet's AutoFilter/_FilterDatabase range to match
# the new (slightly smaller) data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("international_supply_side")
$ws.Activate() | Out-Null

# --- New data rows -------------------------------------------------------
# Columns: A Region, B Medium, C Drive, D Fuel, E New_fuel, F Date, G Reference, H Target
# Rows 12-21 repeat the existing 2030 "ship"/"air" bunker rows (rows 2-11)
# for 2040, rows 22-31 for 2050 and rows 32-41 for 2070.

$rows = @(
    @("all","ship","ship_lpg",      "07_09_lpg",               "16_01_biogas",          2040,0.1,0.3),
    @("all","air", "air_lpg",       "07_09_lpg",               "16_01_biogas",          2040,0.1,0.3),
    @("all","ship","ship_gasoline", "07_01_motor_gasoline",    "16_05_biogasoline",     2040,0.1,0.3),
    @("all","air", "air_gasoline",  "07_01_motor_gasoline",    "16_05_biogasoline",     2040,0.1,0.3),
    @("all","ship","ship_diesel",   "07_07_gas_diesel_oil",    "16_06_biodiesel",       2040,0.1,0.3),
    @("all","air", "air_diesel",    "07_07_gas_diesel_oil",    "16_06_biodiesel",       2040,0.1,0.3),
    @("all","ship","ship_kerosene", "07_06_kerosene",          "16_07_bio_jet_kerosene",2040,0.1,0.3),
    @("all","air", "air_av_gas",    "07_02_aviation_gasoline", "16_07_bio_jet_kerosene",2040,0.1,0.3),
    @("all","air", "air_jet_fuel",  "07_x_jet_fuel",           "16_07_bio_jet_kerosene",2040,0.1,0.3),
    @("all","air", "air_kerosene",  "07_06_kerosene",          "16_07_bio_jet_kerosene",2040,0.1,0.3),

    @("all","ship","ship_lpg",      "07_09_lpg",               "16_01_biogas",          2050,0.1,0.4),
    @("all","air", "air_lpg",       "07_09_lpg",               "16_01_biogas",          2050,0.1,0.4),
    @("all","ship","ship_gasoline", "07_01_motor_gasoline",    "16_05_biogasoline",     2050,0.1,0.4),
    @("all","air", "air_gasoline",  "07_01_motor_gasoline",    "16_05_biogasoline",     2050,0.1,0.4),
    @("all","ship","ship_diesel",   "07_07_gas_diesel_oil",    "16_06_biodiesel",       2050,0.1,0.4),
    @("all","air", "air_diesel",    "07_07_gas_diesel_oil",    "16_06_biodiesel",       2050,0.1,0.4),
    @("all","ship","ship_kerosene", "07_06_kerosene",          "16_07_bio_jet_kerosene",2050,0.1,0.4),
    @("all","air", "air_av_gas",    "07_02_aviation_gasoline", "16_07_bio_jet_kerosene",2050,0.1,0.4),
    @("all","air", "air_jet_fuel",  "07_x_jet_fuel",           "16_07_bio_jet_kerosene",2050,0.1,0.4),
    @("all","air", "air_kerosene",  "07_06_kerosene",          "16_07_bio_jet_kerosene",2050,0.1,0.4),

    @("all","ship","ship_lpg",      "07_09_lpg",               "16_01_biogas",          2070,0.1,0.4),
    @("all","air", "air_lpg",       "07_09_lpg",               "16_01_biogas",          2070,0.1,0.4),
    @("all","ship","ship_gasoline", "07_01_motor_gasoline",    "16_05_biogasoline",     2070,0.1,0.4),
    @("all","air", "air_gasoline",  "07_01_motor_gasoline",    "16_05_biogasoline",     2070,0.1,0.4),
    @("all","ship","ship_diesel",   "07_07_gas_diesel_oil",    "16_06_biodiesel",       2070,0.1,0.4),
    @("all","air", "air_diesel",    "07_07_gas_diesel_oil",    "16_06_biodiesel",       2070,0.1,0.4),
    @("all","ship","ship_kerosene", "07_06_kerosene",          "16_07_bio_jet_kerosene",2070,0.1,0.4),
    @("all","air", "air_av_gas",    "07_02_aviation_gasoline", "16_07_bio_jet_kerosene",2070,0.1,0.4),
    @("all","air", "air_jet_fuel",  "07_x_jet_fuel",           "16_07_bio_jet_kerosene",2070,0.1,0.4),
    @("all","air", "air_kerosene",  "07_06_kerosene",          "16_07_bio_jet_kerosene",2070,0.1,0.4)
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# --- Refresh the autofilter / sort-state over the new (smaller) extent ---
# The filter range is pre-sized well beyond the populated rows (as it was
# before the edit, at row 350) and is only trimmed to 347 here, not to the
# actual last populated row (41).
$lastRow = 347
$dataRange = "A1:H$lastRow"

$ws.AutoFilterMode = $false
$ws.Range($dataRange).AutoFilter() | Out-Null

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("F1:F$lastRow")) | Out-Null
$ws.Sort.SetRange($ws.Range($dataRange))
$ws.Sort.Header = 1
$ws.Sort.Apply() | Out-Null

# --- Shrink the workbook-level _FilterDatabase defined name to match -----
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "international_supply_side!_FilterDatabase") {
        $nm.RefersTo = "=international_supply_side!`$A`$1:`$H`$$lastRow"
    }
}

# --- Restore the selection to where the user left off ---------------------
$ws.Range("J25").Select() | Out-Null

Write-Output "Added $($rows.Count) bunker rows; filter range now $dataRange"
